$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily price record was inserted at the top of the data table
# (row 2), pushing every existing record down by one row. The sheet's
# used range therefore grows from A1:R86 to A1:R87.
$ws.Rows.Item(2).Insert()

# Insert() copies formatting from the row above (the bold header row);
# strip that back out so the new row matches the plain data-row look of
# the rest of the table.
$ws.Rows.Item(2).ClearFormats()

$ws.Cells.Item(2, 1).Value = 3
$ws.Cells.Item(2, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(2, 3).Value = "Coquimbo"
$ws.Cells.Item(2, 4).Value = 44812
$ws.Cells.Item(2, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(2, 5).Value = 5
$ws.Cells.Item(2, 6).Value = 100112035
$ws.Cells.Item(2, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(2, 8).Value = "Sin especificar"
$ws.Cells.Item(2, 9).Value = "Primera"
$ws.Cells.Item(2, 10).Value = 45
$ws.Cells.Item(2, 11).Value = 16000
$ws.Cells.Item(2, 12).Value = 16000
$ws.Cells.Item(2, 13).Value = 16000
$ws.Cells.Item(2, 14).Value = "$/malla 15 kilos"
$ws.Cells.Item(2, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(2, 16).Value = 1067
$ws.Cells.Item(2, 17).Value = 15
$ws.Cells.Item(2, 18).Value = "Hortaliza"
